# Ltbp1-Itgb5.xlsx -- refresh NATMI ligand/receptor expression + specificity
# columns with new TPM-derived numbers (commit: "update scripts wuth new tpm").
# Ligand cols G/H (avg/total expr) vary by sending cluster (col A);
# Receptor cols M/N (avg/total expr) vary by target cluster (col D);
# I/J/O/P/Q/R/S/T are derived specificities/weights recomputed from the above.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("G2").Value = 4.617566333333333
$ws.Range("H2").Value = 13.852699
$ws.Range("I2").Value = 0.07434710147987492
$ws.Range("J2").Value = 0.0743471014798749
$ws.Range("M2").Value = 4.277890333333334
$ws.Range("N2").Value = 12.833671
$ws.Range("O2").Value = 0.04123357425337639
$ws.Range("P2").Value = 0.04123357425337638
$ws.Range("Q2").Value = 19.75344238089212
$ws.Range("R2").Value = 177.780981428029
$ws.Range("S2").Value = 0.003065596729393732
$ws.Range("T2").Value = 0.003065596729393731
# row 3
$ws.Range("G3").Value = 4.617566333333333
$ws.Range("H3").Value = 13.852699
$ws.Range("I3").Value = 0.07434710147987492
$ws.Range("J3").Value = 0.0743471014798749
$ws.Range("O3").Value = 0.4451428460610328
$ws.Range("P3").Value = 0.4451428460610327
$ws.Range("Q3").Value = 213.2510634877335
$ws.Range("R3").Value = 1919.259571389601
$ws.Range("S3").Value = 0.03309508034913995
$ws.Range("T3").Value = 0.03309508034913993
# row 4
$ws.Range("G4").Value = 4.617566333333333
$ws.Range("H4").Value = 13.852699
$ws.Range("I4").Value = 0.07434710147987492
$ws.Range("J4").Value = 0.0743471014798749
$ws.Range("M4").Value = 8.558147333333332
$ws.Range("N4").Value = 25.674442
$ws.Range("O4").Value = 0.08248996024761777
$ws.Range("P4").Value = 0.08248996024761777
$ws.Range("Q4").Value = 39.51781300210644
$ws.Range("R4").Value = 355.660317018958
$ws.Range("S4").Value = 0.006132889445600486
$ws.Range("T4").Value = 0.006132889445600485
# row 5
$ws.Range("G5").Value = 4.617566333333333
$ws.Range("H5").Value = 13.852699
$ws.Range("I5").Value = 0.07434710147987492
$ws.Range("J5").Value = 0.0743471014798749
$ws.Range("M5").Value = 44.72914066666667
$ws.Range("N5").Value = 134.187422
$ws.Range("O5").Value = 0.4311336194379731
$ws.Range("P5").Value = 0.431133619437973
$ws.Range("Q5").Value = 206.5397740613309
$ws.Range("R5").Value = 1858.857966551978
$ws.Range("S5").Value = 0.03205353495574076
$ws.Range("T5").Value = 0.03205353495574075
# row 6
$ws.Range("I6").Value = 0.4103032508824865
$ws.Range("J6").Value = 0.4103032508824864
$ws.Range("M6").Value = 4.277890333333334
$ws.Range("N6").Value = 12.833671
$ws.Range("O6").Value = 0.04123357425337639
$ws.Range("P6").Value = 0.04123357425337638
$ws.Range("Q6").Value = 109.0143591837786
$ws.Range("R6").Value = 981.1292326540072
$ws.Range("S6").Value = 0.01691826956166473
$ws.Range("T6").Value = 0.01691826956166472
# row 7
$ws.Range("I7").Value = 0.4103032508824865
$ws.Range("J7").Value = 0.4103032508824864
$ws.Range("O7").Value = 0.4451428460610328
$ws.Range("P7").Value = 0.4451428460610327
$ws.Range("S7").Value = 0.182643556845924
$ws.Range("T7").Value = 0.1826435568459239
# row 8
$ws.Range("I8").Value = 0.4103032508824865
$ws.Range("J8").Value = 0.4103032508824864
$ws.Range("M8").Value = 8.558147333333332
$ws.Range("N8").Value = 25.674442
$ws.Range("O8").Value = 0.08248996024761777
$ws.Range("P8").Value = 0.08248996024761777
$ws.Range("Q8").Value = 218.0890286209682
$ws.Range("R8").Value = 1962.801257588714
$ws.Range("S8").Value = 0.03384589885476465
$ws.Range("T8").Value = 0.03384589885476465
# row 9
$ws.Range("I9").Value = 0.4103032508824865
$ws.Range("J9").Value = 0.4103032508824864
$ws.Range("M9").Value = 44.72914066666667
$ws.Range("N9").Value = 134.187422
$ws.Range("O9").Value = 0.4311336194379731
$ws.Range("P9").Value = 0.431133619437973
$ws.Range("Q9").Value = 1139.84189090193
$ws.Range("R9").Value = 10258.57701811737
$ws.Range("S9").Value = 0.1768955256201331
$ws.Range("T9").Value = 0.1768955256201331
# row 10
$ws.Range("G10").Value = 31.94872733333333
$ws.Range("H10").Value = 95.846182
$ws.Range("I10").Value = 0.5144041474959183
$ws.Range("J10").Value = 0.5144041474959181
$ws.Range("M10").Value = 4.277890333333334
$ws.Range("N10").Value = 12.833671
$ws.Range("O10").Value = 0.04123357425337639
$ws.Range("P10").Value = 0.04123357425337638
$ws.Range("Q10").Value = 136.6731518215691
$ws.Range("R10").Value = 1230.058366394122
$ws.Range("S10").Value = 0.02121072161201773
$ws.Range("T10").Value = 0.02121072161201772
# row 11
$ws.Range("G11").Value = 31.94872733333333
$ws.Range("H11").Value = 95.846182
$ws.Range("I11").Value = 0.5144041474959183
$ws.Range("J11").Value = 0.5144041474959181
$ws.Range("O11").Value = 0.4451428460610328
$ws.Range("P11").Value = 0.4451428460610327
$ws.Range("Q11").Value = 1475.474219337246
$ws.Range("R11").Value = 13279.26797403522
$ws.Range("S11").Value = 0.2289833262419324
$ws.Range("T11").Value = 0.2289833262419322
# row 12
$ws.Range("G12").Value = 31.94872733333333
$ws.Range("H12").Value = 95.846182
$ws.Range("I12").Value = 0.5144041474959183
$ws.Range("J12").Value = 0.5144041474959181
$ws.Range("M12").Value = 8.558147333333332
$ws.Range("N12").Value = 25.674442
$ws.Range("O12").Value = 0.08248996024761777
$ws.Range("P12").Value = 0.08248996024761777
$ws.Range("Q12").Value = 273.4219156311604
$ws.Range("R12").Value = 2460.797240680444
$ws.Range("S12").Value = 0.04243317767814801
$ws.Range("T12").Value = 0.042433177678148
# row 13
$ws.Range("G13").Value = 31.94872733333333
$ws.Range("H13").Value = 95.846182
$ws.Range("I13").Value = 0.5144041474959183
$ws.Range("J13").Value = 0.5144041474959181
$ws.Range("M13").Value = 44.72914066666667
$ws.Range("N13").Value = 134.187422
$ws.Range("O13").Value = 0.4311336194379731
$ws.Range("P13").Value = 0.431133619437973
$ws.Range("Q13").Value = 1429.039119013645
$ws.Range("R13").Value = 12861.3520711228
$ws.Range("S13").Value = 0.2217769219638202
$ws.Range("T13").Value = 0.2217769219638201
# row 14
$ws.Range("G14").Value = 0.05872333333333333
$ws.Range("H14").Value = 0.17617
$ws.Range("I14").Value = 0.0009455001417203652
$ws.Range("J14").Value = 0.000945500141720365
$ws.Range("M14").Value = 4.277890333333334
$ws.Range("N14").Value = 12.833671
$ws.Range("O14").Value = 0.04123357425337639
$ws.Range("P14").Value = 0.04123357425337638
$ws.Range("Q14").Value = 0.2512119800077778
$ws.Range("R14").Value = 2.26090782007
$ws.Range("S14").Value = 0.00003898635030020458
$ws.Range("T14").Value = 0.00003898635030020456
# row 15
$ws.Range("G15").Value = 0.05872333333333333
$ws.Range("H15").Value = 0.17617
$ws.Range("I15").Value = 0.0009455001417203652
$ws.Range("J15").Value = 0.000945500141720365
$ws.Range("O15").Value = 0.4451428460610328
$ws.Range("P15").Value = 0.4451428460610327
$ws.Range("Q15").Value = 2.711994236981111
$ws.Range("R15").Value = 24.40794813283
$ws.Range("S15").Value = 0.0004208826240365132
$ws.Range("T15").Value = 0.000420882624036513
# row 16
$ws.Range("G16").Value = 0.05872333333333333
$ws.Range("H16").Value = 0.17617
$ws.Range("I16").Value = 0.0009455001417203652
$ws.Range("J16").Value = 0.000945500141720365
$ws.Range("M16").Value = 8.558147333333332
$ws.Range("N16").Value = 25.674442
$ws.Range("O16").Value = 0.08248996024761777
$ws.Range("P16").Value = 0.08248996024761777
$ws.Range("Q16").Value = 0.502562938571111
$ws.Range("R16").Value = 4.52306644714
$ws.Range("S16").Value = 0.0000779942691046299
$ws.Range("T16").Value = 0.00007799426910462988
# row 17
$ws.Range("G17").Value = 0.05872333333333333
$ws.Range("H17").Value = 0.17617
$ws.Range("I17").Value = 0.0009455001417203652
$ws.Range("J17").Value = 0.000945500141720365
$ws.Range("M17").Value = 44.72914066666667
$ws.Range("N17").Value = 134.187422
$ws.Range("O17").Value = 0.4311336194379731
$ws.Range("P17").Value = 0.431133619437973
$ws.Range("Q17").Value = 2.626644237082222
$ws.Range("R17").Value = 23.63979813374
$ws.Range("S17").Value = 0.0004076368982790175
$ws.Range("T17").Value = 0.0004076368982790174
